$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 22 (the "CLO AAA ETF F1" row), pushing the
# existing rows 22-26 down to 23-27.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new "CLO F1 / TRUPS MEZZ F1" entry.
$ws.Cells.Item(22,1).Value2 = "CLO F1"
$ws.Cells.Item(22,2).Value2 = "TRUPS MEZZ F1"
$ws.Cells.Item(22,3).Value2 = 0.07
$ws.Cells.Item(22,4).Value2 = 0.07
$ws.Cells.Item(22,5).Value2 = 0.07

# Match the shaded "detail row" formatting used by the row below (A/B columns).
$ws.Range("A23:B23").Copy()
$ws.Range("A22:B22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection, matching the saved state of the workbook.
$ws.Range("B23").Select() | Out-Null
